$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: update title text
$ws.Range("A17").Value = "Definition of terms"

# Row 19: move exactMatch value out of H19, put closeMatch URL into I19
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = "http://purl.org/aspect/EnviromentalConditionTerms"

# Row 20: set exactMatch URL for wind_speed
$ws.Range("H20").Value = "http://purl.org/aspect/wind_speed"
